$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F (dSF) for specific rows, per the data repull / mean calculation fix.
$updates = @{
    6  = 5
    11 = -3
    14 = -5
    15 = -2
    16 = 3
    17 = 1
    19 = 2
    21 = 6
    23 = 1
    26 = 0
    35 = -1
    39 = -1
    40 = -2
    41 = -7
    42 = -2
    49 = -5
    61 = -7
    62 = -4
    65 = -1
    66 = -3
    67 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
